$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2:E2").ClearContents()
